$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column D ("Service") to make room for "Account number".
# This shifts the existing D:L columns to E:M.
$ws.Range("D1:D3").EntireColumn.Insert()

# Fill in the new "Account number" column (now column D).
$ws.Range("D1").Value = "Account number"
$ws.Range("D2").Value = "{d.meter[i].accountNumber}"
$ws.Range("D3").Value = "{d.meter[i + 1].accountNumber}"
